$wb = $excel.ActiveWorkbook

# --- Update status text: "Ready for handoff" -> "In Translation" ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"
$ws1.Range("E3").Value = "In Translation"
$ws1.Range("F3").Value = "In Translation"

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "In Translation"
$wsZh.Range("C3").Value = "In Translation"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "In Translation"
$wsDe.Range("C3").Value = "In Translation"

# --- Narrow the status columns to fit the shorter text ---
$ws1.Columns("E:F").ColumnWidth = 13.4101845877511
$wsZh.Columns("C:C").ColumnWidth = 13.4101845877511
$wsDe.Columns("C:C").ColumnWidth = 13.4101845877511
